# Update column G ("K" - strikeouts) values on Sheet1 with regenerated data.
# These values come from re-pulling the pitcher's per-game strikeout totals
# (replacing the previous "Strike#" pitch-count-derived values) as part of
# a save_data regen pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 10
    3  = 8
    4  = 3
    5  = 5
    6  = 10
    7  = 3
    8  = 5
    9  = 3
    10 = 0
    11 = 5
    12 = 0
    13 = 4
    14 = 3
    15 = 5
    16 = 7
    17 = 9
    18 = 3
    19 = 5
    20 = 9
    21 = 10
    22 = 4
    23 = 2
    24 = 1
    25 = 6
    26 = 4
    27 = 1
    28 = 2
    29 = 1
    30 = 4
    31 = 0
    32 = 2
    33 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
